$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 1279.4166
$ws.Range("I82").Value = 1279.4166
$ws.Range("K82").Value = 3838.2498
$ws.Range("M82").Value = -3432.2498

$ws.Range("H85").Value = 1279.4166
$ws.Range("I85").Value = 1279.4166
$ws.Range("K85").Value = 3838.2498
$ws.Range("M85").Value = -2434.2498

$ws.Range("H107").Value = 1042.579
$ws.Range("I107").Value = 1001.2857
$ws.Range("J107").Value = 1158.2
$ws.Range("K107").Value = 1001.2857
$ws.Range("L107").Value = 1158.2
$ws.Range("M107").Value = 918.7143
$ws.Range("N107").Value = -4998.2

$ws.Range("H132").Value = 1412.035
$ws.Range("I132").Value = 1191.4906
$ws.Range("J132").Value = 4334.25
$ws.Range("K132").Value = 3574.4718
$ws.Range("L132").Value = 13002.75
$ws.Range("M132").Value = -1044.4718
$ws.Range("N132").Value = -18062.75

$ws.Range("H138").Value = 2654.182
$ws.Range("I138").Value = 2421.7778
$ws.Range("J138").Value = 3700
$ws.Range("K138").Value = 7265.3334
$ws.Range("L138").Value = 11100
$ws.Range("M138").Value = -2125.3334
$ws.Range("N138").Value = -21380

$ws.Range("H139").Value = 68971.5
$ws.Range("J139").Value = 68971.5
$ws.Range("L139").Value = 68971.5
$ws.Range("N139").Value = -79251.5

$ws.Range("H141").Value = 3538.4167
$ws.Range("I141").Value = 3676.2
$ws.Range("K141").Value = 11028.6
$ws.Range("M141").Value = -5848.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4188.2563
$ws.Range("I32").Value = 1754.629
$ws.Range("K32").Value = 1754.629
$ws.Range("M32").Value = -1467.629

$ws.Range("H74").Value = 4950.75
$ws.Range("I74").Value = 2518.7646
$ws.Range("J74").Value = 10857
$ws.Range("K74").Value = 2518.7646
$ws.Range("L74").Value = 10857
$ws.Range("M74").Value = -1644.7646
$ws.Range("N74").Value = -12605

$ws.Range("H77").Value = 4950.75
$ws.Range("I77").Value = 2518.7646
$ws.Range("J77").Value = 10857
$ws.Range("K77").Value = 12593.823
$ws.Range("L77").Value = 54285
$ws.Range("M77").Value = -8225.823
$ws.Range("N77").Value = -63021

$ws.Range("H117").Value = 49111.25
$ws.Range("J117").Value = 49111.25
$ws.Range("L117").Value = 49111.25
$ws.Range("N117").Value = -58289.25

$ws.Range("H127").Value = 90797.336
$ws.Range("J127").Value = 90797.336
$ws.Range("L127").Value = 90797.336
$ws.Range("N127").Value = -100717.336

$ws.Range("H132").Value = 3000.8965
$ws.Range("I132").Value = 2506.875
$ws.Range("J132").Value = 5372.2
$ws.Range("K132").Value = 7520.625
$ws.Range("L132").Value = 16116.6
$ws.Range("M132").Value = -4990.625
$ws.Range("N132").Value = -21176.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 83737.164
$ws.Range("J22").Value = 399.33334
$ws.Range("L22").Value = 399.33334
$ws.Range("N22").Value = -745.33334

$ws.Range("H94").Value = 1344.9615
$ws.Range("I94").Value = 1249.909
$ws.Range("J94").Value = 1867.75
$ws.Range("K94").Value = 1249.909
$ws.Range("L94").Value = 1867.75
$ws.Range("M94").Value = -798.9090000000001
$ws.Range("N94").Value = -2769.75

$ws.Range("H105").Value = 87494.664
$ws.Range("I105").Value = 129773.5
$ws.Range("J105").Value = 2937
$ws.Range("K105").Value = 129773.5
$ws.Range("L105").Value = 2937
$ws.Range("M105").Value = -128026.5
$ws.Range("N105").Value = -6431

$ws.Range("H108").Value = 99896
$ws.Range("J108").Value = 99896
$ws.Range("L108").Value = 99896
$ws.Range("N108").Value = -107576

$ws.Range("H134").Value = 4492.125
$ws.Range("I134").Value = 3287.6
$ws.Range("J134").Value = 6499.6665
$ws.Range("K134").Value = 9862.799999999999
$ws.Range("L134").Value = 19498.9995
$ws.Range("M134").Value = -7327.799999999999
$ws.Range("N134").Value = -24568.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 9431.909
$ws.Range("I7").Value = 9418.091
$ws.Range("J7").Value = 9445.727999999999
$ws.Range("K7").Value = 9418.091
$ws.Range("L7").Value = 9445.727999999999
$ws.Range("M7").Value = -9305.091
$ws.Range("N7").Value = -9671.727999999999

$ws.Range("H31").Value = 1931.2041
$ws.Range("I31").Value = 1298.45
$ws.Range("J31").Value = 2367.5862
$ws.Range("K31").Value = 1298.45
$ws.Range("L31").Value = 2367.5862
$ws.Range("M31").Value = -1003.45
$ws.Range("N31").Value = -2957.5862

$ws.Range("H34").Value = 1931.2041
$ws.Range("I34").Value = 1298.45
$ws.Range("J34").Value = 2367.5862
$ws.Range("K34").Value = 1298.45
$ws.Range("L34").Value = 2367.5862
$ws.Range("M34").Value = -1096.45
$ws.Range("N34").Value = -2771.5862

$ws.Range("H58").Value = 2463
$ws.Range("I58").Value = 1950.9
$ws.Range("J58").Value = 3316.5
$ws.Range("K58").Value = 1950.9
$ws.Range("L58").Value = 3316.5
$ws.Range("M58").Value = -1747.9
$ws.Range("N58").Value = -3722.5

$ws.Range("H86").Value = 1998663.5
$ws.Range("I86").Value = 3582598.5
$ws.Range("J86").Value = 18744.75
$ws.Range("K86").Value = 3582598.5
$ws.Range("L86").Value = 18744.75
$ws.Range("M86").Value = -3581475.5
$ws.Range("N86").Value = -20990.75

$ws.Range("H89").Value = 1998663.5
$ws.Range("I89").Value = 3582598.5
$ws.Range("J89").Value = 18744.75
$ws.Range("K89").Value = 17912992.5
$ws.Range("L89").Value = 93723.75
$ws.Range("M89").Value = -17907376.5
$ws.Range("N89").Value = -104955.75

$ws.Range("H134").Value = 5176547.5
$ws.Range("I134").Value = 11908894
$ws.Range("J134").Value = 127287.375
$ws.Range("K134").Value = 35726682
$ws.Range("L134").Value = 381862.125
$ws.Range("M134").Value = -35724147
$ws.Range("N134").Value = -386932.125

$ws.Range("H136").Value = 2463
$ws.Range("I136").Value = 1950.9
$ws.Range("J136").Value = 3316.5
$ws.Range("K136").Value = 5852.700000000001
$ws.Range("L136").Value = 9949.5
$ws.Range("M136").Value = -3302.700000000001
$ws.Range("N136").Value = -15049.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2215.0833
$ws.Range("J32").Value = 2800
$ws.Range("L32").Value = 8400
$ws.Range("N32").Value = -8966

$ws.Range("H46").Value = 437.5
$ws.Range("J46").Value = 500
$ws.Range("L46").Value = 1500
$ws.Range("N46").Value = -1682

$ws.Range("H107").Value = 1119.0588
$ws.Range("J107").Value = 1168.2858
$ws.Range("L107").Value = 3504.8574
$ws.Range("N107").Value = -7344.857400000001

$ws.Range("H117").Value = 1338.6666
$ws.Range("J117").Value = 2162.5
$ws.Range("L117").Value = 6487.5
$ws.Range("N117").Value = -13371.5

$ws.Range("H131").Value = 1386.8948
$ws.Range("I131").Value = 799.4
$ws.Range("J131").Value = 2039.6666
$ws.Range("K131").Value = 2398.2
$ws.Range("L131").Value = 6118.9998
$ws.Range("M131").Value = 2641.8
$ws.Range("N131").Value = -16198.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 99880.57000000001
$ws.Range("J110").Value = 99880.57000000001
$ws.Range("L110").Value = 99880.57000000001
$ws.Range("N110").Value = -108060.57

$ws.Range("H119").Value = 79316.17999999999
$ws.Range("J119").Value = 79316.17999999999
$ws.Range("L119").Value = 79316.17999999999
$ws.Range("N119").Value = -88992.17999999999

$ws.Range("H122").Value = 5826.7
$ws.Range("I122").Value = 8002.4
$ws.Range("K122").Value = 24007.2
$ws.Range("M122").Value = -21557.2

$ws.Range("H132").Value = 6663.567
$ws.Range("I132").Value = 5631.2856
$ws.Range("J132").Value = 9072.223
$ws.Range("K132").Value = 16893.8568
$ws.Range("L132").Value = 27216.669
$ws.Range("M132").Value = -14363.8568
$ws.Range("N132").Value = -32276.669

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6753.12
$ws.Range("I7").Value = 51214.5
$ws.Range("J7").Value = 2886.913
$ws.Range("K7").Value = 51214.5
$ws.Range("L7").Value = 2886.913
$ws.Range("M7").Value = -51102.5
$ws.Range("N7").Value = -3110.913

$ws.Range("H36").Value = 70000
$ws.Range("J36").Value = 70000
$ws.Range("L36").Value = 70000
$ws.Range("N36").Value = -71124

$ws.Range("H93").Value = 3105.6924
$ws.Range("I93").Value = 2967.1667
$ws.Range("J93").Value = 3224.4285
$ws.Range("K93").Value = 2967.1667
$ws.Range("L93").Value = 3224.4285
$ws.Range("M93").Value = -1719.1667
$ws.Range("N93").Value = -5720.4285

$ws.Range("H126").Value = 6753.12
$ws.Range("I126").Value = 51214.5
$ws.Range("J126").Value = 2886.913
$ws.Range("K126").Value = 153643.5
$ws.Range("L126").Value = 8660.739
$ws.Range("M126").Value = -151173.5
$ws.Range("N126").Value = -13600.739

$ws.Range("H127").Value = 75184.42999999999
$ws.Range("J127").Value = 75184.42999999999
$ws.Range("L127").Value = 75184.42999999999
$ws.Range("N127").Value = -85104.42999999999

$ws.Range("H132").Value = 2782.6155
$ws.Range("I132").Value = 2130.4443
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 6391.3329
$ws.Range("L132").Value = 12750
$ws.Range("N132").Value = -17810
$ws.Range("M132").Value = -3861.3329

$ws.Range("H136").Value = 4489.885
$ws.Range("I136").Value = 4579.9375
$ws.Range("J136").Value = 4345.8
$ws.Range("K136").Value = 13739.8125
$ws.Range("L136").Value = 13037.4
$ws.Range("M136").Value = -11189.8125
$ws.Range("N136").Value = -18137.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3177.577
$ws.Range("I132").Value = 3127.2632
$ws.Range("J132").Value = 3314.1428
$ws.Range("K132").Value = 9381.7896
$ws.Range("L132").Value = 9942.428400000001
$ws.Range("M132").Value = -6851.7896
$ws.Range("N132").Value = -15002.4284
